$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle0 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.344.14"
$ws.Range("D2").Style = $origStyle0
$origStyle1 = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.20%  "
$ws.Range("E2").Style = $origStyle1
$origStyle2 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.011.97"
$ws.Range("D3").Style = $origStyle2
$origStyle3 = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.35%  "
$ws.Range("E3").Style = $origStyle3
$origStyle4 = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E4").Style = $origStyle4
$origStyle5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.88"
$ws.Range("D5").Style = $origStyle5
$origStyle6 = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("E5").Style = $origStyle6
$origStyle7 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.47"
$ws.Range("D6").Style = $origStyle7
$origStyle8 = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.40%  "
$ws.Range("E6").Style = $origStyle8
$origStyle9 = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E7").Style = $origStyle9
$origStyle10 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.006.37"
$ws.Range("D8").Style = $origStyle10
$origStyle11 = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.50%  "
$ws.Range("E8").Style = $origStyle11
$origStyle12 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").Style = $origStyle12
$origStyle13 = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E9").Style = $origStyle13
$origStyle14 = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.95%  "
$ws.Range("E10").Style = $origStyle14
$origStyle15 = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.97%  "
$ws.Range("E11").Style = $origStyle15
$origStyle16 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = $origStyle16
$origStyle17 = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("E12").Style = $origStyle17
$origStyle18 = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.93%  "
$ws.Range("E13").Style = $origStyle18
$origStyle19 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.68"
$ws.Range("D14").Style = $origStyle19
$origStyle20 = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.23%  "
$ws.Range("E14").Style = $origStyle20
$origStyle21 = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E15").Style = $origStyle21
$origStyle22 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.508.39"
$ws.Range("D16").Style = $origStyle22
$origStyle23 = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.38%  "
$ws.Range("E16").Style = $origStyle23
$origStyle24 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.003.27"
$ws.Range("D17").Style = $origStyle24
$origStyle25 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.277.34"
$ws.Range("D18").Style = $origStyle25
$origStyle26 = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.31%  "
$ws.Range("E18").Style = $origStyle26
$origStyle27 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = $origStyle27
$origStyle28 = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("E19").Style = $origStyle28
$origStyle29 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.52"
$ws.Range("D20").Style = $origStyle29
$origStyle30 = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.04%  "
$ws.Range("E20").Style = $origStyle30
$origStyle31 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("D21").Style = $origStyle31
$origStyle32 = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.02%  "
$ws.Range("E21").Style = $origStyle32
$origStyle33 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.671"
$ws.Range("D22").Style = $origStyle33
$origStyle34 = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("E22").Style = $origStyle34
$origStyle35 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.08"
$ws.Range("D23").Style = $origStyle35
$origStyle36 = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.98%  "
$ws.Range("E23").Style = $origStyle36
$origStyle37 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.90"
$ws.Range("D24").Style = $origStyle37
$origStyle38 = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("E24").Style = $origStyle38
$origStyle39 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.46"
$ws.Range("D25").Style = $origStyle39
$origStyle40 = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.69%  "
$ws.Range("E25").Style = $origStyle40
$origStyle41 = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E26").Style = $origStyle41
$origStyle42 = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E27").Style = $origStyle42
$origStyle43 = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.80%  "
$ws.Range("E28").Style = $origStyle43
$origStyle44 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = $origStyle44
$origStyle45 = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("E29").Style = $origStyle45
$origStyle46 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.28"
$ws.Range("D30").Style = $origStyle46
$origStyle47 = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.80%  "
$ws.Range("E30").Style = $origStyle47
$origStyle48 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").Style = $origStyle48
$origStyle49 = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -10.71%  "
$ws.Range("E31").Style = $origStyle49
$origStyle50 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.35"
$ws.Range("D32").Style = $origStyle50
$origStyle51 = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.90%  "
$ws.Range("E32").Style = $origStyle51
$origStyle52 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0948"
$ws.Range("D33").Style = $origStyle52
$origStyle53 = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.41%  "
$ws.Range("E33").Style = $origStyle53
$origStyle54 = $ws.Range("B34").Style
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Mantle"
$ws.Range("B34").Style = $origStyle54
$origStyle55 = $ws.Range("C34").Style
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C34").Style = $origStyle55
$origStyle56 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.940"
$ws.Range("D34").Style = $origStyle56
$origStyle57 = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -8.76%  "
$ws.Range("E34").Style = $origStyle57
$origStyle58 = $ws.Range("B35").Style
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Filecoin"
$ws.Range("B35").Style = $origStyle58
$origStyle59 = $ws.Range("C35").Style
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C35").Style = $origStyle59
$origStyle60 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.61"
$ws.Range("D35").Style = $origStyle60
$origStyle61 = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("E35").Style = $origStyle61
$origStyle62 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.32"
$ws.Range("D36").Style = $origStyle62
$origStyle63 = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("E36").Style = $origStyle63
$origStyle64 = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -14.94%  "
$ws.Range("E37").Style = $origStyle64
$origStyle65 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0676"
$ws.Range("D38").Style = $origStyle65
$origStyle66 = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.37%  "
$ws.Range("E38").Style = $origStyle66
$origStyle67 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.53"
$ws.Range("D39").Style = $origStyle67
$origStyle68 = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.55%  "
$ws.Range("E39").Style = $origStyle68
$origStyle69 = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -9.01%  "
$ws.Range("E40").Style = $origStyle69
$origStyle70 = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("E41").Style = $origStyle70
$origStyle71 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "376.07"
$ws.Range("D42").Style = $origStyle71
$origStyle72 = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.52%  "
$ws.Range("E42").Style = $origStyle72
$origStyle73 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.681.45"
$ws.Range("D43").Style = $origStyle73
$origStyle74 = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.27%  "
$ws.Range("E43").Style = $origStyle74
$origStyle75 = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.46%  "
$ws.Range("E44").Style = $origStyle75
$origStyle76 = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E45").Style = $origStyle76
$origStyle77 = $ws.Range("B46").Style
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Monero"
$ws.Range("B46").Style = $origStyle77
$origStyle78 = $ws.Range("C46").Style
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C46").Style = $origStyle78
$origStyle79 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.35"
$ws.Range("D46").Style = $origStyle79
$origStyle80 = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("E46").Style = $origStyle80
$origStyle81 = $ws.Range("B47").Style
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "TheGraph"
$ws.Range("B47").Style = $origStyle81
$origStyle82 = $ws.Range("C47").Style
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C47").Style = $origStyle82
$origStyle83 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.235"
$ws.Range("D47").Style = $origStyle83
$origStyle84 = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.01%  "
$ws.Range("E47").Style = $origStyle84
$origStyle85 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("D48").Style = $origStyle85
$origStyle86 = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.77%  "
$ws.Range("E48").Style = $origStyle86
$origStyle87 = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("E49").Style = $origStyle87
$origStyle88 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.53"
$ws.Range("D50").Style = $origStyle88
$origStyle89 = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.73%  "
$ws.Range("E50").Style = $origStyle89
$origStyle90 = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.91%  "
$ws.Range("E51").Style = $origStyle90
